# Update the "Out of PO" roster table (A2:C18) with the refreshed player list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("Donovan Mitchell",   "PG,SG",    "Cleveland Cavaliers"),
  @("Dyson Daniels",      "PG,SG,SF", "Atlanta Hawks"),
  @("Jaden McDaniels",    "SF,PF",    "Minnesota Timberwolves"),
  @("De'Andre Hunter",    "SF,PF",    "Atlanta Hawks"),
  @("Alperen Sengün",     "C",        "Houston Rockets"),
  @("Jonas Valanciunas",  "C",        "Washington Wizards"),
  @("Domantas Sabonis",   "C",        "Sacramento Kings"),
  @("Michael Porter Jr.", "SF,PF",    "Denver Nuggets"),
  @("Kristaps Porzingis", "PF,C",     "Boston Celtics"),
  @("Kel'el Ware",        "PF,C",     "Miami Heat"),
  @("Josh Hart",          "SG,SF,PF", "New York Knicks"),
  @("Kelly Oubre Jr.",    "SG,SF",    "Philadelphia 76ers"),
  @("Victor Wembanyama",  "C",        "San Antonio Spurs"),
  @("Andrew Wiggins",     "SF,PF",    "Golden State Warriors"),
  @("Malik Beasley",      "SG,SF",    "Detroit Pistons"),
  @("Cam Thomas",         "SG,SF",    "Brooklyn Nets"),
  @("Donte DiVincenzo",   "PG,SG,SF", "Washington Wizards")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
